# Refresh the crypto price/volume table for the 29-12-2022 07:xx snapshot.
# Numeric-looking "Price" cells (column D) are written through Formula with a
# leading apostrophe (then restyled "Normal") so they stay plain text cells
# (matching the original inlineStr storage) instead of being auto-coerced to
# floating point numbers by the COM Value setter.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = "'244.74"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(3, 4).Formula = "'23.94"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(4, 4).Formula = "'5.201"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(5, 4).Formula = "'0.05735"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Formula = "'6.492"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(7, 4).Formula = "'3.168"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(8, 4).Formula = "'0.8136"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(10, 2).Value = 'WazirX'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(10, 4).Formula = "'0.1371"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '9WazirXWRX'
$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(11, 4).Formula = "'0.06952"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '10MandalaExchangeTokenMDX'
$ws.Cells.Item(12, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(12, 4).Formula = "'0.03179"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Cells.Item(13, 2).Value = 'BitrueCoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(13, 4).Formula = "'0.02928"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '12BitrueCoinBTR'
$ws.Cells.Item(14, 2).Value = 'BitMartToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(14, 4).Formula = "'0.09326"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '13BitMartTokenBMX'
$ws.Cells.Item(15, 2).Value = 'MCDex'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Cells.Item(15, 4).Formula = "'3.846"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '14MCDexMCB'
$ws.Cells.Item(16, 2).Value = 'BitForexToken'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(16, 4).Formula = "'0.001531"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '15BitForexTokenBF'
$ws.Cells.Item(17, 2).Value = 'CoinExToken'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Cells.Item(17, 4).Formula = "'0.04695"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '16CoinExTokenCET'
$ws.Cells.Item(18, 2).Value = 'One'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(18, 4).Formula = "'0.0006015"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '17OneONE'
$ws.Cells.Item(19, 4).Formula = "'0.006188"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(20, 4).Formula = "'0.001244"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(21, 4).Formula = "'0.004105"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(22, 2).Value = 'NitroEx'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Cells.Item(22, 4).Formula = "'0.00007002"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '21NitroExNTX'
$ws.Cells.Item(23, 2).Value = 'LEO'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(23, 4).Formula = "'3.551"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '22LEOLEO'
$ws.Cells.Item(24, 2).Value = 'BTSEToken'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(24, 4).Formula = "'2.151"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '23BTSETokenBTSE'
$ws.Cells.Item(25, 2).Value = 'BitpandaEcosystemToken'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Cells.Item(25, 4).Formula = "'0.3193"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '24BitpandaEcosystemTokenBEST'
$ws.Cells.Item(26, 2).Value = 'ProBitToken'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Cells.Item(26, 4).Formula = "'0.1330"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '25ProBitTokenPROB'
$ws.Cells.Item(27, 4).Formula = "'0.0002330"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(40, 4).Formula = "'0.03710"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(41, 2).Value = 'BKEXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Cells.Item(41, 4).Formula = "'0.1052"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '40BKEXTokenBKK'
$ws.Cells.Item(42, 2).Value = 'CEJI'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Cells.Item(42, 4).Formula = "'0.002290"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '41CEJICEJI'
$ws.Cells.Item(43, 2).Value = 'KickToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Cells.Item(43, 4).Formula = "'0.003062"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '42KickTokenKICKWorstin24h'
$ws.Cells.Item(44, 4).Formula = "'0.008083"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(46, 4).Formula = "'0.00000000751"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 4).Formula = "'0.4543"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(48, 4).Formula = "'0.002627"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '47BOLOBOLOBestin24h'
$ws.Cells.Item(49, 4).Formula = "'0.00002102"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(50, 4).Formula = "'0.0002002"
$ws.Cells.Item(50, 4).Style = "Normal"
